$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells stay text (avoid Excel auto-numeric coercion) while updating values
$priceCells = @("D2", "D3", "D5", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D27", "D31", "D32", "D33", "D34", "D36", "D37", "D40", "D41", "D43", "D46", "D47", "D50", "D51")
foreach ($ref in $priceCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "41.141.66"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.138.40"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "234.64"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("E6").Value = "  -5.01%  "
$ws.Range("D7").Value = "68.97"
$ws.Range("E7").Value = "  -6.06%  "
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -7.42%  "
$ws.Range("D10").Value = "38.25"
$ws.Range("E10").Value = "  -10.78%  "
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  -7.33%  "
$ws.Range("D12").Value = "53.11"
$ws.Range("E12").Value = "  -7.66%  "
$ws.Range("D13").Value = "0.0992"
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").Value = "6.54"
$ws.Range("E14").Value = "  -7.36%  "
$ws.Range("D15").Value = "2.458.88"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").Value = "14.29"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "2.129.82"
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("D18").Value = "0.772"
$ws.Range("E18").Value = "  -7.96%  "
$ws.Range("D19").Value = "41.000.88"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  -7.99%  "
$ws.Range("D21").Value = "68.65"
$ws.Range("E21").Value = "  -5.61%  "
$ws.Range("D22").Value = "5.68"
$ws.Range("E22").Value = "  -8.27%  "
$ws.Range("D23").Value = "223.74"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("E24").Value = "  -13.04%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -9.69%  "
$ws.Range("D27").Value = "10.48"
$ws.Range("E27").Value = "  -11.13%  "
$ws.Range("E28").Value = "  -9.30%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -6.77%  "
$ws.Range("D31").Value = "169.17"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "19.48"
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("D33").Value = "30.63"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "0.0745"
$ws.Range("E34").Value = "  -6.77%  "
$ws.Range("E35").Value = "  -12.48%  "
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "0.0997"
$ws.Range("E37").Value = "  -9.49%  "
$ws.Range("E38").Value = "  -5.41%  "
$ws.Range("E39").Value = "  -8.29%  "
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  -4.92%  "
$ws.Range("D41").Value = "11.56"
$ws.Range("E41").Value = "  -17.70%  "
$ws.Range("E42").Value = "  -7.95%  "
$ws.Range("D43").Value = "56.87"
$ws.Range("E43").Value = "  -13.21%  "
$ws.Range("E44").Value = "  -7.33%  "
$ws.Range("E45").Value = "  -8.37%  "
$ws.Range("D46").Value = "0.0947"
$ws.Range("E46").Value = "  -5.91%  "
$ws.Range("D47").Value = "96.35"
$ws.Range("E47").Value = "  -8.44%  "
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("E49").Value = "  -6.46%  "
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "2.12"
$ws.Range("E51").Value = "  -12.07%  "

foreach ($ref in $priceCells) { $ws.Range($ref).Style = "Normal" }
